$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh coin prices / 1h volume percentages and re-order a handful of
# rows, matching the latest scrape from coinranking.com (GitHub Actions).

$ws.Range('D2').Value = "'60.176.49"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.79%  '
$ws.Range('D3').Value = "'3.310.16"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.86%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = "'560.17"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.24%  '
$ws.Range('D6').Value = "'144.73"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.39%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = "'3.311.64"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.85%  '
$ws.Range('D9').Value = "'0.484"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.18%  '
$ws.Range('E10').Value = '  -2.62%  '
$ws.Range('E11').Value = '  -2.52%  '
$ws.Range('E12').Value = '  -0.88%  '
$ws.Range('D13').Value = "'3.873.63"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.96%  '
$ws.Range('E14').Value = '  +0.81%  '
$ws.Range('D15').Value = "'27.34"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.34%  '
$ws.Range('D16').Value = "'3.332.43"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.15%  '
$ws.Range('D17').Value = "'0.0000167"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.63%  '
$ws.Range('D18').Value = "'60.157.27"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.89%  '
$ws.Range('D19').Value = "'6.18"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.26%  '
$ws.Range('D20').Value = "'14.34"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.15%  '
$ws.Range('D21').Value = "'8.70"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.33%  '
$ws.Range('D22').Value = "'375.11"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.84%  '
$ws.Range('D23').Value = "'74.28"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.10%  '
$ws.Range('D24').Value = "'0.552"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.41%  '
$ws.Range('D25').Value = "'1.00"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('D26').Value = "'3.482.56"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.79%  '
$ws.Range('D27').Value = "'0.0000106"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -7.48%  '
$ws.Range('D28').Value = "'0.172"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.74%  '
$ws.Range('D30').Value = "'7.26"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.91%  '
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('D33').Value = "'7.64"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.52%  '
$ws.Range('D34').Value = "'22.66"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.27%  '
$ws.Range('D35').Value = "'1.30"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.15%  '
$ws.Range('E36').Value = '  -3.47%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = "'1.54"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.18%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').Value = "'166.77"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.34%  '
$ws.Range('D39').Value = "'6.78"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.14%  '
$ws.Range('D40').Value = "'27.91"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -12.35%  '
$ws.Range('D41').Value = "'3.337.39"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.04%  '
$ws.Range('D42').Value = "'0.0739"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.30%  '
$ws.Range('D43').Value = "'41.92"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.55%  '
$ws.Range('D44').Value = "'0.753"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.90%  '
$ws.Range('D45').Value = "'4.23"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.44%  '
$ws.Range('B46').Value = 'ONDO'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D46').Value = "'1.13"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.68%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').Value = "'1.60"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.85%  '
$ws.Range('D48').Value = "'2.398.25"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.78%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').Value = "'6.62"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.91%  '
$ws.Range('B50').Value = 'FirstDigitalUSD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D50').Value = "'0.998"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.20%  '
$ws.Range('D51').Value = "'21.91"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.02%  '
